$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Test1"
$ws.Range("A2").Value = "Test2"
$ws.Range("B1").Value = "TestB1"
$ws.Range("A3").Value = "Test3"
$ws.Range("A4").Value = "Test4"
$ws.Range("A5").Value = "Test5"
$ws.Range("B2").Value = "TestB2"
$ws.Range("B4").Value = "TestB4"
$ws.Range("B5").Value = "TestB5"
$ws.Range("B3").Value = "TestB3"

$ws.Range("A1:XFD1").Select()
